function Add-FileHyperlink {
    param($ws, $cellRef, $displayText, $targetUrl)
    $ws.Range($cellRef).Value = $displayText
    $ws.Hyperlinks.Add($ws.Range($cellRef), $targetUrl, [System.Type]::Missing, [System.Type]::Missing, $displayText)
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = 0xED9564
}

$wb = $excel.ActiveWorkbook
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# --- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("B3").Value = "Handed back: in sync with en-US"
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("B3").Value = "Handed back: in sync with en-US"

# --- zh-cn: populate "Latest Target File" (E) / "Latest Handback File" (F) ---
Add-FileHyperlink $zh "E2" "98634d4f-a790-4ca3-b6e6-59723ededfab.md" "https://github.com/OpenLocalizationTest/oltest/blob/943f5880a1a27cea6366acb0c4bc7c1e1051cd19/e2e/98634d4f-a790-4ca3-b6e6-59723ededfab.md"
Add-FileHyperlink $zh "F2" "98634d4f-a790-4ca3-b6e6-59723ededfab.f600a82c0fb23017c042922fb9979d1b198fff34.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90c2f5d1b51e32f59843f0dbceae57d3c89bf8b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/98634d4f-a790-4ca3-b6e6-59723ededfab.f600a82c0fb23017c042922fb9979d1b198fff34.zh-cn.xlf"
Add-FileHyperlink $zh "E3" "9915040f-7ef0-4485-8a38-b65208a2b685.md" "https://github.com/OpenLocalizationTest/oltest/blob/943f5880a1a27cea6366acb0c4bc7c1e1051cd19/e2e/9915040f-7ef0-4485-8a38-b65208a2b685.md"
Add-FileHyperlink $zh "F3" "9915040f-7ef0-4485-8a38-b65208a2b685.366bd0d6ad1f02ea7e248ef2200971494bc1e2fd.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/90c2f5d1b51e32f59843f0dbceae57d3c89bf8b7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/9915040f-7ef0-4485-8a38-b65208a2b685.366bd0d6ad1f02ea7e248ef2200971494bc1e2fd.zh-cn.xlf"

# --- zh-cn: "Latest Handback DateTime" (G) now has a real timestamp ---
$zh.Range("G2").Value = "2016-03-10 06:12:42"
$zh.Range("G3").Value = "2016-03-10 06:12:42"

# --- de-de: populate "Latest Target File" (E) / "Latest Handback File" (F) ---
Add-FileHyperlink $de "E2" "98634d4f-a790-4ca3-b6e6-59723ededfab.md" "https://github.com/OpenLocalizationTest/oltest/blob/943f5880a1a27cea6366acb0c4bc7c1e1051cd19/e2e/98634d4f-a790-4ca3-b6e6-59723ededfab.md"
Add-FileHyperlink $de "F2" "98634d4f-a790-4ca3-b6e6-59723ededfab.f600a82c0fb23017c042922fb9979d1b198fff34.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/043b3d7bf132250ab5c3f6c8dc6dc99426b96399/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/98634d4f-a790-4ca3-b6e6-59723ededfab.f600a82c0fb23017c042922fb9979d1b198fff34.de-de.xlf"
Add-FileHyperlink $de "E3" "9915040f-7ef0-4485-8a38-b65208a2b685.md" "https://github.com/OpenLocalizationTest/oltest/blob/943f5880a1a27cea6366acb0c4bc7c1e1051cd19/e2e/9915040f-7ef0-4485-8a38-b65208a2b685.md"
Add-FileHyperlink $de "F3" "9915040f-7ef0-4485-8a38-b65208a2b685.366bd0d6ad1f02ea7e248ef2200971494bc1e2fd.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/043b3d7bf132250ab5c3f6c8dc6dc99426b96399/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high/9915040f-7ef0-4485-8a38-b65208a2b685.366bd0d6ad1f02ea7e248ef2200971494bc1e2fd.de-de.xlf"

# --- de-de: "Latest Handback DateTime" (G) now has a real timestamp ---
$de.Range("G2").Value = "2016-03-10 06:12:51"
$de.Range("G3").Value = "2016-03-10 06:12:51"
